$wb = $excel.ActiveWorkbook

# ----- Sheet 1: Summary -----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6170411985018727
$ws1.Range("C2").Value = 0.5724217844727694
$ws1.Range("D2").Value = 0.9250936329588015
$ws1.Range("E2").Value = 0.7072297780959198
$ws1.Range("F2").Value = 0.8236078692897633
$ws1.Range("G2").Value = 0.9036797298248083
$ws1.Range("H2").Value = 0.7728559104490174
$ws1.Range("I2").Value = 494
$ws1.Range("J2").Value = 369
$ws1.Range("K2").Value = 165
$ws1.Range("L2").Value = 40

# ----- Sheet 2: Classification Report -----
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.8048780487804879
$ws2.Range("C2").Value = 0.3089887640449438
$ws2.Range("D2").Value = 0.4465493910690122

$ws2.Range("B3").Value = 0.5724217844727694
$ws2.Range("C3").Value = 0.9250936329588015
$ws2.Range("D3").Value = 0.7072297780959198

$ws2.Range("B4").Value = 0.6170411985018727
$ws2.Range("C4").Value = 0.6170411985018727
$ws2.Range("D4").Value = 0.6170411985018727
$ws2.Range("E4").Value = 0.6170411985018727

$ws2.Range("B5").Value = 0.6886499166266287
$ws2.Range("C5").Value = 0.6170411985018727
$ws2.Range("D5").Value = 0.576889584582466

$ws2.Range("B6").Value = 0.6886499166266286
$ws2.Range("C6").Value = 0.6170411985018727
$ws2.Range("D6").Value = 0.5768895845824661

# ----- Sheet 3: Confusion Matrix -----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 165
$ws3.Range("C2").Value = 369
$ws3.Range("B3").Value = 40
$ws3.Range("C3").Value = 494
